$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"

$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("B3").Value = "Ready for handoff"
$ws2.Range("D3").Value = "2016-03-09 16:11:11"

$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("B3").Value = "Ready for handoff"
$ws3.Range("D3").Value = "2016-03-09 16:11:14"
